# Rename the three "Reference*" sheets to the shorter "Ref*" names.
# Excel automatically rewrites every formula that referenced the old sheet
# names (e.g. ReferenceMetadata!B4 -> RefMetadata!B4) as part of the rename,
# exactly like renaming a sheet tab interactively in the UI.
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("ReferenceProperties").Name = "RefProperties"
$wb.Worksheets.Item("ReferenceClasses").Name = "RefClasses"
$wb.Worksheets.Item("ReferenceMetadata").Name = "RefMetadata"

# The workbook was last left open on the renamed "RefMetadata" sheet
# (previously "ReferenceProperties" was the selected tab) -- make that the
# active tab/sheet, moving the tabSelected flag accordingly.
$wb.Worksheets.Item("RefMetadata").Activate()
